$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 9796.762000000001
$ws.Range("I15").Value = 9796.762000000001
$ws.Range("K15").Value = 29390.286
$ws.Range("M15").Value = -29221.286
$ws.Range("H33").Value = 405.2857
$ws.Range("J33").Value = 1309
$ws.Range("L33").Value = 1309
$ws.Range("N33").Value = -1767
$ws.Range("H51").Value = 8944.056
$ws.Range("J51").Value = 9213.857
$ws.Range("L51").Value = 9213.857
$ws.Range("N51").Value = -10181.857
$ws.Range("H57").Value = 33888.332
$ws.Range("I57").Value = 70000
$ws.Range("J57").Value = 29374.375
$ws.Range("K57").Value = 210000
$ws.Range("L57").Value = 88123.125
$ws.Range("M57").Value = -209501
$ws.Range("N57").Value = -89121.125
$ws.Range("H74").Value = 9604.166999999999
$ws.Range("I74").Value = 4875
$ws.Range("K74").Value = 4875
$ws.Range("M74").Value = -3939
$ws.Range("H76").Value = 13500
$ws.Range("J76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15630
$ws.Range("H77").Value = 9604.166999999999
$ws.Range("I77").Value = 4875
$ws.Range("K77").Value = 24375
$ws.Range("M77").Value = -19695
$ws.Range("H79").Value = 13500
$ws.Range("J79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17184
$ws.Range("H80").Value = 535.4286
$ws.Range("I80").Value = 524.6667
$ws.Range("K80").Value = 1574.0001
$ws.Range("M80").Value = -576.0001
$ws.Range("H83").Value = 535.4286
$ws.Range("I83").Value = 524.6667
$ws.Range("K83").Value = 4722.0003
$ws.Range("M83").Value = 269.9997000000003
$ws.Range("H92").Value = 176.92308
$ws.Range("I92").Value = 176.92308
$ws.Range("K92").Value = 176.92308
$ws.Range("M92").Value = 1071.07692
$ws.Range("H107").Value = 1442.8
$ws.Range("I107").Value = 1492
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1492
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 428
$ws.Range("N107").Value = -4840
$ws.Range("H115").Value = 2798.5
$ws.Range("I115").Value = 597
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 1791
$ws.Range("L115").Value = 15000
$ws.Range("M115").Value = -224
$ws.Range("N115").Value = -18134
$ws.Range("H137").Value = 1697.25
$ws.Range("I137").Value = 1263
$ws.Range("K137").Value = 3789
$ws.Range("M137").Value = -1239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 7013
$ws.Range("I33").Value = 4026
$ws.Range("K33").Value = 4026
$ws.Range("M33").Value = -3697
$ws.Range("H45").Value = 7088.636
$ws.Range("I45").Value = 3868.75
$ws.Range("J45").Value = 8928.571
$ws.Range("K45").Value = 3868.75
$ws.Range("L45").Value = 8928.571
$ws.Range("M45").Value = -3491.75
$ws.Range("N45").Value = -9682.571
$ws.Range("H63").Value = 5767.4287
$ws.Range("I63").Value = 3174.75
$ws.Range("K63").Value = 3174.75
$ws.Range("M63").Value = -2488.75
$ws.Range("H66").Value = 5767.4287
$ws.Range("I66").Value = 3174.75
$ws.Range("K66").Value = 15873.75
$ws.Range("M66").Value = -12441.75
$ws.Range("H74").Value = 3374.7334
$ws.Range("I74").Value = 2973.862
$ws.Range("K74").Value = 2973.862
$ws.Range("M74").Value = -2099.862
$ws.Range("H77").Value = 3374.7334
$ws.Range("I77").Value = 2973.862
$ws.Range("K77").Value = 14869.31
$ws.Range("M77").Value = -10501.31

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 320
$ws.Range("I11").Value = 31.666666
$ws.Range("K11").Value = 31.666666
$ws.Range("M11").Value = 108.333334
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H102").Value = 4551
$ws.Range("I102").Value = 4551
$ws.Range("K102").Value = 4551
$ws.Range("M102").Value = -1306
$ws.Range("H134").Value = 3435.158
$ws.Range("I134").Value = 3496.4688
$ws.Range("J134").Value = 3108.1667
$ws.Range("K134").Value = 10489.4064
$ws.Range("L134").Value = 9324.500100000001
$ws.Range("M134").Value = -7954.4064
$ws.Range("N134").Value = -14394.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("K6").Value = 5000
$ws.Range("M6").Value = -4887
$ws.Range("H13").Value = 2500
$ws.Range("I13").Value = 2500
$ws.Range("K13").Value = 2500
$ws.Range("M13").Value = -2361
$ws.Range("H31").Value = 4237.353
$ws.Range("I31").Value = 4270.857
$ws.Range("J31").Value = 4213.9
$ws.Range("K31").Value = 4270.857
$ws.Range("L31").Value = 4213.9
$ws.Range("M31").Value = -3975.857
$ws.Range("N31").Value = -4803.9
$ws.Range("H34").Value = 4237.353
$ws.Range("I34").Value = 4270.857
$ws.Range("J34").Value = 4213.9
$ws.Range("K34").Value = 4270.857
$ws.Range("L34").Value = 4213.9
$ws.Range("M34").Value = -4068.857
$ws.Range("N34").Value = -4617.9
$ws.Range("H92").Value = 64866.668
$ws.Range("J92").Value = 64866.668
$ws.Range("L92").Value = 64866.668
$ws.Range("N92").Value = -69858.66800000001
$ws.Range("H94").Value = 2364.6667
$ws.Range("I94").Value = 1547.5
$ws.Range("J94").Value = 3999
$ws.Range("K94").Value = 1547.5
$ws.Range("L94").Value = 3999
$ws.Range("M94").Value = -1096.5
$ws.Range("N94").Value = -4901
$ws.Range("H132").Value = 3014.2307
$ws.Range("I132").Value = 3098.75
$ws.Range("K132").Value = 9296.25
$ws.Range("M132").Value = -6766.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9788.125
$ws.Range("I3").Value = 8309.23
$ws.Range("J3").Value = 16196.667
$ws.Range("K3").Value = 24927.69
$ws.Range("L3").Value = 48590.001
$ws.Range("M3").Value = -24815.69
$ws.Range("N3").Value = -48814.001
$ws.Range("H128").Value = 593873.75
$ws.Range("I128").Value = 593873.75
$ws.Range("K128").Value = 1781621.25
$ws.Range("M128").Value = -1776641.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3825.1765
$ws.Range("I97").Value = 787.7857
$ws.Range("J97").Value = 17999.666
$ws.Range("K97").Value = 787.7857
$ws.Range("L97").Value = 17999.666
$ws.Range("M97").Value = -291.7857
$ws.Range("N97").Value = -18991.666
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 60121.94
$ws.Range("I61").Value = 60121.94
$ws.Range("K61").Value = 60121.94
$ws.Range("M61").Value = -59919.94
$ws.Range("H81").Value = 29999
$ws.Range("I81").Value = 29999
$ws.Range("K81").Value = 29999
$ws.Range("M81").Value = -29001
$ws.Range("H84").Value = 29999
$ws.Range("I84").Value = 29999
$ws.Range("K84").Value = 89997
$ws.Range("M84").Value = -85005
$ws.Range("H113").Value = 60121.94
$ws.Range("I113").Value = 60121.94
$ws.Range("K113").Value = 60121.94
$ws.Range("M113").Value = -57951.94
$ws.Range("H132").Value = 8414.5
$ws.Range("I132").Value = 8912.1
$ws.Range("K132").Value = 26736.3
$ws.Range("M132").Value = -24206.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 87440.11
$ws.Range("J5").Value = 87440.11
$ws.Range("L5").Value = 87440.11
$ws.Range("N5").Value = -87664.11
$ws.Range("H107").Value = 659.4375
$ws.Range("I107").Value = 566.7692
$ws.Range("J107").Value = 1061
$ws.Range("K107").Value = 1700.3076
$ws.Range("L107").Value = 3183
$ws.Range("M107").Value = 219.6924000000001
$ws.Range("N107").Value = -7023
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530
$ws.Range("H132").Value = 2673
$ws.Range("I132").Value = 2673
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8019
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5489
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 5550.3887
$ws.Range("J136").Value = 6200
$ws.Range("L136").Value = 18600
$ws.Range("N136").Value = -23700
